$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell-level updates per the diff (row-by-row) for the cryptos table
$ws.Range("D2").Value = "96.542.74"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "3.674.18"
$ws.Range("E3").Value = "  +2.48%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.44"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.84"
$ws.Range("E6").Value = "  +10.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "661.35"
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.423"
$ws.Range("E8").Value = "  +2.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.08"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "3.671.73"
$ws.Range("E11").Value = "  +2.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.40"
$ws.Range("E12").Value = "  +4.17%  "
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.89"
$ws.Range("E14").Value = "  +6.89%  "
$ws.Range("D15").Value = "4.356.93"
$ws.Range("E15").Value = "  +2.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000269"
$ws.Range("E16").Value = "  +4.64%  "
$ws.Range("D17").Value = "96.412.45"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.91"
$ws.Range("E18").Value = "  +14.79%  "
$ws.Range("D19").Value = "3.663.36"
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.02"
$ws.Range("E20").Value = "  +2.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.35"
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.525"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "530.07"
$ws.Range("E23").Value = "  +3.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.43"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000204"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.95"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.58"
$ws.Range("E27").Value = "  +4.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.02"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.166"
$ws.Range("E29").Value = "  +10.88%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.53"
$ws.Range("E30").Value = "  +8.84%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.05"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("B33").Value = "Cronos"
$ws.Range("C33").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.186"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  +13.89%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.06"
$ws.Range("E35").Value = "  +5.05%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "635.53"
$ws.Range("E37").Value = "  +1.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.590"
$ws.Range("E38").Value = "  +3.96%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "45.78"
$ws.Range("E39").Value = "  +38.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.76"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.161"
$ws.Range("E41").Value = "  +5.11%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.968"
$ws.Range("E42").Value = "  +6.18%  "
$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.97"
$ws.Range("E43").Value = "  +5.30%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.33"
$ws.Range("E44").Value = "  +8.06%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0456"
$ws.Range("E46").Value = "  +6.24%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.436"
$ws.Range("E47").Value = "  +16.72%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.29"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.63"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("B50").Value = "MantraDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.64"
$ws.Range("E50").Value = "  +3.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.59"
$ws.Range("E51").Value = "  +2.76%  "
